$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the two new rows, identical to the existing rows 2-6.
$colA = "products__item\ in-stock\ products__item_3-in-row`"]:nth-child(2) [type=`"button"
$colB = "New!iPhone 15 PlusFrom:`$ 1,200.00`$"
$colC = "Increased"
$colD = " "
$colE = "\31 52173-case-650"
$colF = "\31 52174-case-655"
$colG = "256 GB"
$colH = "Green"
$colI = "Apple"
$colJ = "'1"

foreach ($r in 7,8) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $colD
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ
}
